$wb = $excel.ActiveWorkbook

# "All" sheet: append a new row for the "tutorial" job, and leave the
# selection on row 2 (full row A2:D2) as recorded in the saved view state.
$wsAll = $wb.Worksheets.Item("All")
$wsAll.Range("A4").Value = "tutorial"
$wsAll.Range("A2:D2").Select() | Out-Null

# "Job to Run" sheet is the tab that is active/selected when the workbook
# is saved, with the cursor on A9.
$wsJob = $wb.Worksheets.Item("Job to Run")
$wsJob.Activate() | Out-Null
$wsJob.Range("A9").Select() | Out-Null
